$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1164.9667
$ws.Range("I40").Value = 1118.04
$ws.Range("J40").Value = 1399.6
$ws.Range("K40").Value = 1118.04
$ws.Range("L40").Value = 1399.6
$ws.Range("M40").Value = -943.04
$ws.Range("N40").Value = -1749.6

$ws.Range("H100").Value = 12822535
$ws.Range("J100").Value = 2993.2222
$ws.Range("L100").Value = 2993.2222
$ws.Range("N100").Value = -4075.2222

$ws.Range("H113").Value = 2914.889
$ws.Range("I113").Value = 2666.8
$ws.Range("J113").Value = 3225
$ws.Range("K113").Value = 2666.8
$ws.Range("L113").Value = 3225
$ws.Range("M113").Value = 587.1999999999998
$ws.Range("N113").Value = -9733

$ws.Range("H137").Value = 715.6604
$ws.Range("I137").Value = 607.67645
$ws.Range("J137").Value = 908.8946999999999
$ws.Range("K137").Value = 1823.02935
$ws.Range("L137").Value = 2726.6841
$ws.Range("M137").Value = 726.9706499999998
$ws.Range("N137").Value = -7826.6841

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 929.2449
$ws.Range("I61").Value = 632.5122
$ws.Range("J61").Value = 2450
$ws.Range("K61").Value = 632.5122
$ws.Range("L61").Value = 2450
$ws.Range("M61").Value = -420.5122
$ws.Range("N61").Value = -2874

$ws.Range("H80").Value = 10071.857
$ws.Range("J80").Value = 16000.75
$ws.Range("L80").Value = 16000.75
$ws.Range("N80").Value = -17996.75

$ws.Range("H83").Value = 10071.857
$ws.Range("J83").Value = 16000.75
$ws.Range("L83").Value = 48002.25
$ws.Range("N83").Value = -57986.25

$ws.Range("H97").Value = 3264
$ws.Range("I97").Value = 3330
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 3330
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -2834
$ws.Range("N97").Value = -3992

$ws.Range("H102").Value = 3657.05
$ws.Range("I102").Value = 2882
$ws.Range("J102").Value = 5982.2
$ws.Range("K102").Value = 2882
$ws.Range("L102").Value = 5982.2
$ws.Range("M102").Value = -1260
$ws.Range("N102").Value = -9226.200000000001

$ws.Range("H122").Value = 32280090
$ws.Range("I122").Value = 34506132
$ws.Range("K122").Value = 103518396
$ws.Range("M122").Value = -103515946

$ws.Range("H136").Value = 929.2449
$ws.Range("I136").Value = 632.5122
$ws.Range("J136").Value = 2450
$ws.Range("K136").Value = 1897.5366
$ws.Range("L136").Value = 7350
$ws.Range("M136").Value = 652.4634000000001
$ws.Range("N136").Value = -12450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 812.3333
$ws.Range("I94").Value = 758.8333
$ws.Range("J94").Value = 1133.3334
$ws.Range("K94").Value = 758.8333
$ws.Range("L94").Value = 1133.3334
$ws.Range("M94").Value = -307.8333
$ws.Range("N94").Value = -2035.3334

$ws.Range("H107").Value = 1007.8182
$ws.Range("I107").Value = 785.5454999999999
$ws.Range("J107").Value = 1452.3636
$ws.Range("K107").Value = 785.5454999999999
$ws.Range("L107").Value = 1452.3636
$ws.Range("M107").Value = 1134.4545
$ws.Range("N107").Value = -5292.3636

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 8600.5
$ws.Range("J50").Value = 8600.5
$ws.Range("L50").Value = 8600.5
$ws.Range("N50").Value = -9850.5

$ws.Range("H74").Value = 17125.8
$ws.Range("J74").Value = 17125.8
$ws.Range("L74").Value = 17125.8
$ws.Range("N74").Value = -18873.8

$ws.Range("H77").Value = 17125.8
$ws.Range("J77").Value = 17125.8
$ws.Range("L77").Value = 51377.39999999999
$ws.Range("N77").Value = -60113.39999999999

$ws.Range("H107").Value = 908.875
$ws.Range("I107").Value = 711.8333
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 711.8333
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 1208.1667
$ws.Range("N107").Value = -5340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 165.76666
$ws.Range("I2").Value = 22.533333
$ws.Range("J2").Value = 309
$ws.Range("K2").Value = 135.199998
$ws.Range("L2").Value = 1854
$ws.Range("M2").Value = -22.19999799999999
$ws.Range("N2").Value = -2080

$ws.Range("H38").Value = 1466.5883
$ws.Range("I38").Value = 930.1818
$ws.Range("J38").Value = 2450
$ws.Range("K38").Value = 2790.5454
$ws.Range("L38").Value = 7350
$ws.Range("M38").Value = -2443.5454
$ws.Range("N38").Value = -8044

$ws.Range("H113").Value = 1017.4821
$ws.Range("I113").Value = 891.2
$ws.Range("J113").Value = 1029.8628
$ws.Range("K113").Value = 2673.6
$ws.Range("L113").Value = 3089.588400000001
$ws.Range("M113").Value = -503.6000000000004
$ws.Range("N113").Value = -7429.588400000001

$ws.Range("H133").Value = 6567.9473
$ws.Range("I133").Value = 4038.4285
$ws.Range("J133").Value = 7139.129
$ws.Range("K133").Value = 12115.2855
$ws.Range("L133").Value = 21417.387
$ws.Range("M133").Value = -7055.2855
$ws.Range("N133").Value = -31537.387

$ws.Range("H137").Value = 5391.325
$ws.Range("I137").Value = 20807
$ws.Range("J137").Value = 3189.0857
$ws.Range("K137").Value = 62421
$ws.Range("L137").Value = 9567.257100000001
$ws.Range("M137").Value = -57321
$ws.Range("N137").Value = -19767.2571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1551.625
$ws.Range("I97").Value = 1036.6666
$ws.Range("J97").Value = 2213.7144
$ws.Range("K97").Value = 1036.6666
$ws.Range("L97").Value = 2213.7144
$ws.Range("M97").Value = -540.6666
$ws.Range("N97").Value = -3205.7144

$ws.Range("H107").Value = 347.17392
$ws.Range("I107").Value = 271.42856
$ws.Range("J107").Value = 465
$ws.Range("K107").Value = 271.42856
$ws.Range("L107").Value = 465
$ws.Range("M107").Value = 1648.57144
$ws.Range("N107").Value = -4305

$ws.Range("H122").Value = 2663.3572
$ws.Range("I122").Value = 2635.875
$ws.Range("J122").Value = 2700
$ws.Range("K122").Value = 7907.625
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -5457.625
$ws.Range("N122").Value = -13000

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 276.6154
$ws.Range("I16").Value = 276.6154
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 276.6154
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -106.6154
$ws.Range("N16").ClearContents()

$ws.Range("H40").Value = 1857.7407
$ws.Range("I40").Value = 1627.4286
$ws.Range("J40").Value = 2105.7693
$ws.Range("K40").Value = 1627.4286
$ws.Range("L40").Value = 2105.7693
$ws.Range("M40").Value = -1491.4286
$ws.Range("N40").Value = -2377.7693

$ws.Range("H93").Value = 1692.3704
$ws.Range("I93").Value = 1629.2941
$ws.Range("J93").Value = 1799.6
$ws.Range("K93").Value = 1629.2941
$ws.Range("L93").Value = 1799.6
$ws.Range("M93").Value = -381.2941000000001
$ws.Range("N93").Value = -4295.6

$ws.Range("H122").Value = 2742.4375
$ws.Range("I122").Value = 3075
$ws.Range("J122").Value = 2631.5833
$ws.Range("K122").Value = 9225
$ws.Range("L122").Value = 7894.749899999999
$ws.Range("M122").Value = -6775
$ws.Range("N122").Value = -12794.7499

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 104140
$ws.Range("I96").Value = 1575
$ws.Range("J96").Value = 172516.67
$ws.Range("K96").Value = 1575
$ws.Range("L96").Value = 172516.67
$ws.Range("M96").Value = -202
$ws.Range("N96").Value = -175262.67

$ws.Range("H107").Value = 734.1667
$ws.Range("I107").Value = 780.4
$ws.Range("J107").Value = 503
$ws.Range("K107").Value = 2341.2
$ws.Range("L107").Value = 1509
$ws.Range("M107").Value = -421.1999999999998
$ws.Range("N107").Value = -5349

$ws.Range("H122").Value = 1929.1428
$ws.Range("I122").Value = 1863.5
$ws.Range("J122").Value = 2016.6666
$ws.Range("K122").Value = 5590.5
$ws.Range("L122").Value = 6049.9998
$ws.Range("M122").Value = -3140.5
$ws.Range("N122").Value = -10949.9998
